{"js": "// The \"Author\" paragraph style is currently bold (w:b in its run\n// properties). The edit removes that bold formatting from the style\n// definition so paragraphs using the \"Author\" style (e.g. the \"Author\"\n// placeholder paragraph) are no longer rendered bold.\nconst styles = context.document.getStyles();\nconst authorStyle = styles.getByNameOrNullObject(\"Author\");\nawait context.sync();\n\nif (authorStyle.isNullObject) {\n  throw new Error('Style \"Author\" was not found in the document.');\n}\n\nauthorStyle.font.bold = false;\nawait context.sync();\n", "ps1": "# The \"Author\" paragraph style is currently bold. Remove the bold\n# formatting from the style definition itself (not a direct/run-level\n# override) so every paragraph that uses the \"Author\" style renders\n# without bold.\n$d = $word.ActiveDocument\n$authorStyle = $d.Styles(\"Author\")\n$authorStyle.Font.Bold = $false\n"}
